# Fruta / hortaliza, semanal
# Insert a new weekly price observation row at row 72 (Pera - Packham's Triumph,
# Primera) on the "Feria Lagunitas de Puerto Montt" sheet. Inserting shifts every
# subsequent row (old 72..111) down by one (new 73..112).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(72).Insert()

$ws.Range("A72").Value = 4
$ws.Range("B72").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C72").Value = "Los Lagos"
$ws.Range("D72").Value = Get-Date -Year 2021 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100104
$ws.Range("H72").Value = "Frutos de pepita"
$ws.Range("I72").Value = 100104005
$ws.Range("J72").Value = "Pera"
$ws.Range("K72").Value = "Packham's Triumph"
$ws.Range("L72").Value = "Primera"
$ws.Range("M72").Value = 200
$ws.Range("N72").Value = 16000
$ws.Range("O72").Value = 16000
$ws.Range("P72").Value = 16000
$ws.Range("Q72").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R72").Value = "Región de O'Higgins"
$ws.Range("S72").Value = 1067
$ws.Range("T72").Value = 15
